$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 49 with the full "Gas Station" entry (problem solved).
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "Array"
$ws.Range("C49").Value = "Gas Station"
$ws.Range("D49").Value = "Medium"
$ws.Range("E49").Value = "Done"
$ws.Range("F49").Value = 45889
$ws.Range("F48").Copy($ws.Range("F49"))
$ws.Range("G49").Value = "O(n)"
$ws.Range("H49").Value = "O(1)"
$ws.Range("I49").Value = "Normal Traversing"

# Remove the stray "Candy" placeholder row entirely (row 50).
$ws.Range("C50").Clear()

# Update the saved view state (scroll position / selection).
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("G44").Select() | Out-Null
